$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Sema5a"
$ws.Range("C2").Value = "Plxnb3"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3843716666666667
$ws.Range("H2").Value = 1.153115
$ws.Range("I2").Value = 0.009283037010184481
$ws.Range("J2").Value = 0.009283037010184483
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.4299933333333333
$ws.Range("N2").Value = 1.28998
$ws.Range("O2").Value = 0.1062132503660503
$ws.Range("P2").Value = 0.1062132503660503
$ws.Range("Q2").Value = 0.1652772541888889
$ws.Range("R2").Value = 1.4874952877
$ws.Range("S2").Value = 0.0009859815341200352
$ws.Range("T2").Value = 0.0009859815341200355

$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Sema5a"
$ws.Range("C3").Value = "Plxnb3"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3843716666666667
$ws.Range("H3").Value = 1.153115
$ws.Range("I3").Value = 0.009283037010184481
$ws.Range("J3").Value = 0.009283037010184483
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.7426423333333334
$ws.Range("N3").Value = 2.227927
$ws.Range("O3").Value = 0.1834411140081888
$ws.Range("P3").Value = 0.1834411140081888
$ws.Range("Q3").Value = 0.2854506714005556
$ws.Range("R3").Value = 2.569056042605001
$ws.Range("S3").Value = 0.001702890650527487
$ws.Range("T3").Value = 0.001702890650527488

$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Sema5a"
$ws.Range("C4").Value = "Plxnb3"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.3843716666666667
$ws.Range("H4").Value = 1.153115
$ws.Range("I4").Value = 0.009283037010184481
$ws.Range("J4").Value = 0.009283037010184483
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 2.875760666666666
$ws.Range("N4").Value = 8.627281999999999
$ws.Range("O4").Value = 0.710345635625761
$ws.Range("P4").Value = 0.710345635625761
$ws.Range("Q4").Value = 1.105360920381111
$ws.Range("R4").Value = 9.948248283430001
$ws.Range("S4").Value = 0.00659416482553696
$ws.Range("T4").Value = 0.00659416482553696

$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Sema5a"
$ws.Range("C5").Value = "Plxnb3"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 26.63881666666667
$ws.Range("H5").Value = 79.91645
$ws.Range("I5").Value = 0.6433593900630531
$ws.Range("J5").Value = 0.6433593900630532
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.4299933333333333
$ws.Range("N5").Value = 1.28998
$ws.Range("O5").Value = 0.1062132503660503
$ws.Range("P5").Value = 0.1062132503660503
$ws.Range("Q5").Value = 11.45451357455555
$ws.Range("R5").Value = 103.090622171
$ws.Range("S5").Value = 0.06833329197211647
$ws.Range("T5").Value = 0.06833329197211647

$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Sema5a"
$ws.Range("C6").Value = "Plxnb3"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 26.63881666666667
$ws.Range("H6").Value = 79.91645
$ws.Range("I6").Value = 0.6433593900630531
$ws.Range("J6").Value = 0.6433593900630532
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 0.7426423333333334
$ws.Range("N6").Value = 2.227927
$ws.Range("O6").Value = 0.1834411140081888
$ws.Range("P6").Value = 0.1834411140081888
$ws.Range("Q6").Value = 19.78311296657223
$ws.Range("R6").Value = 178.04801669915
$ws.Range("S6").Value = 0.1180185632207953
$ws.Range("T6").Value = 0.1180185632207953

$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Sema5a"
$ws.Range("C7").Value = "Plxnb3"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 26.63881666666667
$ws.Range("H7").Value = 79.91645
$ws.Range("I7").Value = 0.6433593900630531
$ws.Range("J7").Value = 0.6433593900630532
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 2.875760666666666
$ws.Range("N7").Value = 8.627281999999999
$ws.Range("O7").Value = 0.710345635625761
$ws.Range("P7").Value = 0.710345635625761
$ws.Range("Q7").Value = 76.60686117654444
$ws.Range("R7").Value = 689.4617505888999
$ws.Range("S7").Value = 0.4570075348701414
$ws.Range("T7").Value = 0.4570075348701415

$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Sema5a"
$ws.Range("C8").Value = "Plxnb3"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 14.38262166666667
$ws.Range("H8").Value = 43.147865
$ws.Range("I8").Value = 0.3473575729267623
$ws.Range("J8").Value = 0.3473575729267623
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.4299933333333333
$ws.Range("N8").Value = 1.28998
$ws.Range("O8").Value = 0.1062132503660503
$ws.Range("P8").Value = 0.1062132503660503
$ws.Range("Q8").Value = 6.184431432522221
$ws.Range("R8").Value = 55.65988289269999
$ws.Range("S8").Value = 0.03689397685981378
$ws.Range("T8").Value = 0.03689397685981378

$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Sema5a"
$ws.Range("C9").Value = "Plxnb3"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 14.38262166666667
$ws.Range("H9").Value = 43.147865
$ws.Range("I9").Value = 0.3473575729267623
$ws.Range("J9").Value = 0.3473575729267623
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 0.7426423333333334
$ws.Range("N9").Value = 2.227927
$ws.Range("O9").Value = 0.1834411140081888
$ws.Range("P9").Value = 0.1834411140081888
$ws.Range("Q9").Value = 10.68114371398389
$ws.Range("R9").Value = 96.130293425855
$ws.Range("S9").Value = 0.06371966013686596
$ws.Range("T9").Value = 0.06371966013686596

$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Sema5a"
$ws.Range("C10").Value = "Plxnb3"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 14.38262166666667
$ws.Range("H10").Value = 43.147865
$ws.Range("I10").Value = 0.3473575729267623
$ws.Range("J10").Value = 0.3473575729267623
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 2.875760666666666
$ws.Range("N10").Value = 8.627281999999999
$ws.Range("O10").Value = 0.710345635625761
$ws.Range("P10").Value = 0.710345635625761
$ws.Range("Q10").Value = 41.36097767254777
$ws.Range("R10").Value = 372.2487990529299
$ws.Range("S10").Value = 0.2467439359300826
$ws.Range("T10").Value = 0.2467439359300826

